$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row 1: rename existing L/M headers, add new N:Q headers ---
$ws.Range("L1").Value = "Estimate lower bound"
$ws.Range("M1").Value = "Estimate upper bound"
$ws.Range("N1").Value = "Initial guess lower bound"
$ws.Range("O1").Value = "Initial guess upper bound"
$ws.Range("P1").Value = "Constraints lower bound"
$ws.Range("Q1").Value = "Constraints upper bound"

# Copy the header style (bold, bordered, centered) from L1 onto the new header cells N1:Q1
$ws.Range("L1").Copy()
$ws.Range("N1:Q1").PasteSpecial(-4122)

# Row 2: Probability of contact within household layer
$ws.Range("A2").Value = "Probability of contact within household layer"
$ws.Range("L2").Value = 0.9813744
$ws.Range("M2").Value = 0.9816373
$ws.Range("N2").Value = 0.9514584
$ws.Range("O2").Value = 0.9997778
$ws.Range("P2").Value = 0.95
$ws.Range("Q2").Value = 1

# Row 3: Consecutive daily contact probability within household layer
$ws.Range("A3").Value = "Consecutive daily contact probability within household layer"
$ws.Range("L3").Value = 0.4035243
$ws.Range("M3").Value = 0.406743
$ws.Range("N3").Value = 0.4018447
$ws.Range("O3").Value = 0.9878015
$ws.Range("P3").Value = 0.4
$ws.Range("Q3").Value = 1

# Row 4: Contact probability when healthy within household layer
$ws.Range("A4").Value = "Contact probability when healthy within household layer"
$ws.Range("L4").Value = 0.0549288
$ws.Range("M4").Value = 0.057753
$ws.Range("N4").Value = 0.0519458
$ws.Range("O4").Value = 0.4999198
$ws.Range("P4").Value = 0.05
$ws.Range("Q4").Value = 0.5

# Row 5: Contact probability when symptomatic within household layer
$ws.Range("A5").Value = "Contact probability when symptomatic within household layer"
$ws.Range("L5").Value = 0.0212009
$ws.Range("M5").Value = 0.0216763
$ws.Range("N5").Value = 0.0108859
$ws.Range("O5").Value = 0.0994009
$ws.Range("P5").Value = 0.01
$ws.Range("Q5").Value = 0.1

# Row 6: Steepness of the logistic contact probability function within household layer
$ws.Range("A6").Value = "Steepness of the logistic contact probability function within household layer"
$ws.Range("L6").Value = 3.7553894
$ws.Range("M6").Value = 3.8905682
$ws.Range("N6").Value = 1.2364132
$ws.Range("O6").Value = 19.9753985
$ws.Range("P6").Value = 1
$ws.Range("Q6").Value = 20

# Row 7: Phase relative to symptom onset for symptomatic (days) within household layer
$ws.Range("A7").Value = "Phase relative to symptom onset for symptomatic (days) within household layer"
$ws.Range("L7").Value = 4.6765422
$ws.Range("M7").Value = 4.7717213
$ws.Range("N7").Value = 0.0986915
$ws.Range("O7").Value = 9.9010677
$ws.Range("P7").Value = 0
$ws.Range("Q7").Value = 10

# Row 8: Phase relative to symptom onset for resuming normal social context (days) within household layer
$ws.Range("A8").Value = "Phase relative to symptom onset for resuming normal social context (days) within household layer"
$ws.Range("L8").Value = 7.1105597
$ws.Range("M8").Value = 7.1653285
$ws.Range("N8").Value = 0.0678561
$ws.Range("O8").Value = 9.8411842
$ws.Range("P8").Value = 0
$ws.Range("Q8").Value = 10

# Row 9: Probability of contact within school layer
$ws.Range("A9").Value = "Probability of contact within school layer"
$ws.Range("L9").Value = 0.8011814
$ws.Range("M9").Value = 0.8064416
$ws.Range("N9").Value = 0.1137086
$ws.Range("O9").Value = 0.9714956
$ws.Range("P9").Value = 0.1
$ws.Range("Q9").Value = 1

# Row 10: Consecutive daily contact probability within school layer
$ws.Range("A10").Value = "Consecutive daily contact probability within school layer"
$ws.Range("L10").Value = 0.9074275
$ws.Range("M10").Value = 0.9098335
$ws.Range("N10").Value = 0.5020798
$ws.Range("O10").Value = 0.9911082
$ws.Range("P10").Value = 0.5
$ws.Range("Q10").Value = 1

# Row 11: Contact probability when healthy within school layer
$ws.Range("A11").Value = "Contact probability when healthy within school layer"
$ws.Range("L11").Value = 0.0127297
$ws.Range("M11").Value = 0.0158094
$ws.Range("N11").Value = 0.0109628
$ws.Range("O11").Value = 0.4999261
$ws.Range("P11").Value = 0.01
$ws.Range("Q11").Value = 0.5

# Row 12: Contact probability when symptomatic within school layer
$ws.Range("A12").Value = "Contact probability when symptomatic within school layer"
$ws.Range("L12").Value = 0.0429207
$ws.Range("M12").Value = 0.0432555
$ws.Range("N12").Value = 0.0018008
$ws.Range("O12").Value = 0.0490501
$ws.Range("P12").Value = 0.001
$ws.Range("Q12").Value = 0.05

# Row 13: Steepness of the logistic contact probability function within school layer
$ws.Range("A13").Value = "Steepness of the logistic contact probability function within school layer"
$ws.Range("L13").Value = 6.7171996
$ws.Range("M13").Value = 6.7834145
$ws.Range("N13").Value = 1.0637532
$ws.Range("O13").Value = 9.4453913
$ws.Range("P13").Value = 1
$ws.Range("Q13").Value = 10

# Row 14: Phase relative to symptom onset for symptomatic (days) within school layer
$ws.Range("A14").Value = "Phase relative to symptom onset for symptomatic (days) within school layer"
$ws.Range("L14").Value = 9.9340981
$ws.Range("M14").Value = 10
$ws.Range("N14").Value = 0.1014381
$ws.Range("O14").Value = 9.9581626
$ws.Range("P14").Value = 0
$ws.Range("Q14").Value = 10

# Row 15: Phase relative to symptom onset for resuming normal social context (days) within school layer
$ws.Range("A15").Value = "Phase relative to symptom onset for resuming normal social context (days) within school layer"
$ws.Range("L15").Value = 1.7065511
$ws.Range("M15").Value = 1.7727431
$ws.Range("N15").Value = 0.0041295
$ws.Range("O15").Value = 9.8999123
$ws.Range("P15").Value = 0
$ws.Range("Q15").Value = 10

# Row 16: Probability of contact within workplace layer
$ws.Range("A16").Value = "Probability of contact within workplace layer"
$ws.Range("L16").Value = 0.2141153
$ws.Range("M16").Value = 0.2188239
$ws.Range("N16").Value = 0.1010826
$ws.Range("O16").Value = 0.9967066
$ws.Range("P16").Value = 0.1
$ws.Range("Q16").Value = 1

# Row 17: Consecutive daily contact probability within workplace layer
$ws.Range("A17").Value = "Consecutive daily contact probability within workplace layer"
$ws.Range("L17").Value = 0.5632769
$ws.Range("M17").Value = 0.5669683
$ws.Range("N17").Value = 0.5039645
$ws.Range("O17").Value = 0.9994351
$ws.Range("P17").Value = 0.5
$ws.Range("Q17").Value = 1

# Row 18: Contact probability when healthy within workplace layer
$ws.Range("A18").Value = "Contact probability when healthy within workplace layer"
$ws.Range("L18").Value = 0.0355397
$ws.Range("M18").Value = 0.0380466
$ws.Range("N18").Value = 0.0122231
$ws.Range("O18").Value = 0.4913164
$ws.Range("P18").Value = 0.01
$ws.Range("Q18").Value = 0.5

# Row 19: Contact probability when symptomatic within workplace layer
$ws.Range("A19").Value = "Contact probability when symptomatic within workplace layer"
$ws.Range("L19").Value = 0.0201447
$ws.Range("M19").Value = 0.0206327
$ws.Range("N19").Value = 0.002221
$ws.Range("O19").Value = 0.049448
$ws.Range("P19").Value = 0.001
$ws.Range("Q19").Value = 0.05

# Row 20: Steepness of the logistic contact probability function within workplace layer
$ws.Range("A20").Value = "Steepness of the logistic contact probability function within workplace layer"
$ws.Range("L20").Value = 2.7617305
$ws.Range("M20").Value = 2.8686074
$ws.Range("N20").Value = 1.072661
$ws.Range("O20").Value = 9.7834415
$ws.Range("P20").Value = 1
$ws.Range("Q20").Value = 10

# Row 21: Phase relative to symptom onset for symptomatic (days) within workplace layer
$ws.Range("A21").Value = "Phase relative to symptom onset for symptomatic (days) within workplace layer"
$ws.Range("L21").Value = 2.6870014
$ws.Range("M21").Value = 2.7917888
$ws.Range("N21").Value = 0.0046317
$ws.Range("O21").Value = 9.9969221
$ws.Range("P21").Value = 0
$ws.Range("Q21").Value = 10

# Row 22: Phase relative to symptom onset for resuming normal social context (days) within workplace layer
$ws.Range("A22").Value = "Phase relative to symptom onset for resuming normal social context (days) within workplace layer"
$ws.Range("L22").Value = 6.8170766
$ws.Range("M22").Value = 6.8714812
$ws.Range("N22").Value = 0.0332962
$ws.Range("O22").Value = 9.8447942
$ws.Range("P22").Value = 0
$ws.Range("Q22").Value = 10

# Row 23: Probability of contact within heath care layer
$ws.Range("A23").Value = "Probability of contact within heath care layer"
$ws.Range("L23").Value = 0.0059791
$ws.Range("M23").Value = 0.0060079
$ws.Range("N23").Value = 0.005001
$ws.Range("O23").Value = 0.0079966
$ws.Range("P23").Value = 0.005
$ws.Range("Q23").Value = 0.008

# Row 24: Consecutive daily contact probability within heath care layer
$ws.Range("A24").Value = "Consecutive daily contact probability within heath care layer"
$ws.Range("L24").Value = 0.8415189
$ws.Range("M24").Value = 0.8473494
$ws.Range("N24").Value = 0.5064355
$ws.Range("O24").Value = 0.9935533
$ws.Range("P24").Value = 0.5
$ws.Range("Q24").Value = 1

# Row 25: Contact probability when healthy within heath care layer
$ws.Range("A25").Value = "Contact probability when healthy within heath care layer"
$ws.Range("L25").Value = 0.0022543
$ws.Range("M25").Value = 0.0023189
$ws.Range("N25").Value = 0.0010241
$ws.Range("O25").Value = 0.0059853
$ws.Range("P25").Value = 0.001
$ws.Range("Q25").Value = 0.006

# Row 26: Contact probability when symptomatic within heath care layer
$ws.Range("A26").Value = "Contact probability when symptomatic within heath care layer"
$ws.Range("L26").Value = 0.7168134
$ws.Range("M26").Value = 0.719144
$ws.Range("N26").Value = 0.6049105
$ws.Range("O26").Value = 0.9980445
$ws.Range("P26").Value = 0.6
$ws.Range("Q26").Value = 1

# Row 27: Steepness of the logistic contact probability function within heath care layer
$ws.Range("A27").Value = "Steepness of the logistic contact probability function within heath care layer"
$ws.Range("L27").Value = 4.1730533
$ws.Range("M27").Value = 4.225386
$ws.Range("N27").Value = 1.1089668
$ws.Range("O27").Value = 9.9162264
$ws.Range("P27").Value = 1
$ws.Range("Q27").Value = 10

# Row 28: Phase relative to symptom onset for symptomatic (days) within heath care layer
$ws.Range("A28").Value = "Phase relative to symptom onset for symptomatic (days) within heath care layer"
$ws.Range("L28").Value = 9.345397
$ws.Range("M28").Value = 9.3850803
$ws.Range("N28").Value = 0.0174672
$ws.Range("O28").Value = 9.9391382
$ws.Range("P28").Value = 0
$ws.Range("Q28").Value = 10

# Row 29: Phase relative to symptom onset for resuming normal social context (days) within heath care layer
$ws.Range("A29").Value = "Phase relative to symptom onset for resuming normal social context (days) within heath care layer"
$ws.Range("L29").Value = 4.3731756
$ws.Range("M29").Value = 4.438816
$ws.Range("N29").Value = 0.1601804
$ws.Range("O29").Value = 9.9270844
$ws.Range("P29").Value = 0
$ws.Range("Q29").Value = 10

# Row 30: Probability of contact within municipality layer
$ws.Range("A30").Value = "Probability of contact within municipality layer"
$ws.Range("L30").Value = 0.0000074
$ws.Range("M30").Value = 0.0000074
$ws.Range("N30").Value = 0.0000021
$ws.Range("O30").Value = 0.0000099
$ws.Range("P30").Value = 0.000002
$ws.Range("Q30").Value = 0.00001

# Row 31: Consecutive daily contact probability within municipality layer
$ws.Range("A31").Value = "Consecutive daily contact probability within municipality layer"
$ws.Range("L31").Value = 0.5677729
$ws.Range("M31").Value = 0.5725763
$ws.Range("N31").Value = 0.5081613
$ws.Range("O31").Value = 0.995623
$ws.Range("P31").Value = 0.5
$ws.Range("Q31").Value = 1

# Row 32: Contact probability when healthy within municipality layer
$ws.Range("A32").Value = "Contact probability when healthy within municipality layer"
$ws.Range("L32").Value = 0.01047
$ws.Range("M32").Value = 0.0108936
$ws.Range("N32").Value = 0.010305
$ws.Range("O32").Value = 0.0997263
$ws.Range("P32").Value = 0.01
$ws.Range("Q32").Value = 0.1

# Row 33: Contact probability when symptomatic within municipality layer
$ws.Range("A33").Value = "Contact probability when symptomatic within municipality layer"
$ws.Range("L33").Value = 0.048963
$ws.Range("M33").Value = 0.0494151
$ws.Range("N33").Value = 0.0011979
$ws.Range("O33").Value = 0.0499194
$ws.Range("P33").Value = 0.001
$ws.Range("Q33").Value = 0.05

# Row 34: Steepness of the logistic contact probability function within municipality layer
$ws.Range("A34").Value = "Steepness of the logistic contact probability function within municipality layer"
$ws.Range("L34").Value = 3.3325891
$ws.Range("M34").Value = 3.400308
$ws.Range("N34").Value = 1.1983502
$ws.Range("O34").Value = 9.8401297
$ws.Range("P34").Value = 1
$ws.Range("Q34").Value = 10

# Row 35: Phase relative to symptom onset for symptomatic (days) within municipality layer
$ws.Range("A35").Value = "Phase relative to symptom onset for symptomatic (days) within municipality layer"
$ws.Range("L35").Value = 7.33048
$ws.Range("M35").Value = 7.4362635
$ws.Range("N35").Value = 0.1208321
$ws.Range("O35").Value = 9.7650716
$ws.Range("P35").Value = 0
$ws.Range("Q35").Value = 10

# Row 36: Phase relative to symptom onset for resuming normal social context (days) within municipality layer
$ws.Range("A36").Value = "Phase relative to symptom onset for resuming normal social context (days) within municipality layer"
$ws.Range("L36").Value = 9.8913702
$ws.Range("M36").Value = 9.9477663
$ws.Range("N36").Value = 0.0159508
$ws.Range("O36").Value = 9.9239655
$ws.Range("P36").Value = 0
$ws.Range("Q36").Value = 10

# Row 37: Overdispersion rate
$ws.Range("A37").Value = "Overdispersion rate"
$ws.Range("L37").Value = 0.166497
$ws.Range("M37").Value = 0.1679448
$ws.Range("N37").Value = 0.0025337
$ws.Range("O37").Value = 0.19394
$ws.Range("P37").Value = 0
$ws.Range("Q37").Value = 0.2

# Row 38: Overdispersion weight
$ws.Range("A38").Value = "Overdispersion weight"
$ws.Range("L38").Value = 6.2801237
$ws.Range("M38").Value = 6.4135898
$ws.Range("N38").Value = 1.1236438
$ws.Range("O38").Value = 19.3580253
$ws.Range("P38").Value = 1
$ws.Range("Q38").Value = 20
